$wb = $excel.ActiveWorkbook

# Rename sheets (task order timestamps updated)
$wb.Worksheets.Item(1).Name = "GNG_TO-16504778444126682"
$wb.Worksheets.Item(2).Name = "NB_TO-16504778463847609"
$wb.Worksheets.Item(3).Name = "RS_TO-16504778463866704"
$wb.Worksheets.Item(4).Name = "TOL_TO-1650477846433669"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16504778464967034"

# Sheet 1 (GNG)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16504778443716695.csv"
$ws1.Range("B3").Value = "GNG_stims-16504778443966699.csv"
$ws1.Range("B4").Value = "go_stims-16504778443976715.csv"
$ws1.Range("B5").Value = "GNG_stims-16504778444117045.csv"

# Sheet 2 (NB)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-1650477845967672.csv"
$ws2.Range("B3").Value = "OB-16504778451827047.csv"
$ws2.Range("B4").Value = "ZB-match_2-16504778449146707.csv"
$ws2.Range("B5").Value = "ZB-match_9-16504778445327053.csv"
$ws2.Range("B6").Value = "TB-1650477846356706.csv"
$ws2.Range("B7").Value = "ZB-match_4-1650477844742668.csv"
$ws2.Range("B8").Value = "TB-16504778461987052.csv"
$ws2.Range("B9").Value = "OB-1650477845397704.csv"
$ws2.Range("B10").Value = "OB-16504778452717035.csv"

# Sheet 3 (RS)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# Sheet 4 (TOL)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-165047784640074.csv"
$ws4.Range("B3").Value = "ZM_stims-165047784638967.csv"
$ws4.Range("B4").Value = "MM_stims-16504778464167068.csv"
$ws4.Range("B5").Value = "ZM_stims-16504778464016712.csv"
$ws4.Range("B6").Value = "MM_stims-16504778464327044.csv"
$ws4.Range("B7").Value = "ZM_stims-16504778464167068.csv"

# Sheet 5 (vSAT)
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16504778464646697.csv"
$ws5.Range("B3").Value = "SAT_stims-16504778464366693.csv"
$ws5.Range("B4").Value = "SAT_stims-16504778464487033.csv"
$ws5.Range("B5").Value = "vSAT_stims-16504778464807055.csv"
